$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 14639554
$ws.Cells.Item(132, 9).Value = 17120458
$ws.Cells.Item(132, 11).Value = 51361374
$ws.Cells.Item(132, 13).Value = -51358844
$ws.Cells.Item(137, 8).Value = 1059.821
$ws.Cells.Item(137, 9).Value = 898.918
$ws.Cells.Item(137, 10).Value = 1348.5
$ws.Cells.Item(137, 11).Value = 2696.754
$ws.Cells.Item(137, 12).Value = 4045.5
$ws.Cells.Item(137, 13).Value = -146.7539999999999
$ws.Cells.Item(137, 14).Value = -9145.5
$ws.Cells.Item(138, 8).Value = 44603.895
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 44603.895
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 133811.685
$ws.Cells.Item(138, 14).Value = -144091.685
$ws.Cells.Item(141, 8).Value = 6862.5
$ws.Cells.Item(141, 9).Value = 6420.5884
$ws.Cells.Item(141, 10).Value = 9366.666999999999
$ws.Cells.Item(141, 11).Value = 19261.7652
$ws.Cells.Item(141, 12).Value = 28100.001
$ws.Cells.Item(141, 13).Value = -14081.7652
$ws.Cells.Item(141, 14).Value = -38460.001
$ws.Cells.Item(138, 13).ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10828.12
$ws.Cells.Item(32, 9).Value = 9325.529
$ws.Cells.Item(32, 10).Value = 19342.8
$ws.Cells.Item(32, 11).Value = 9325.529
$ws.Cells.Item(32, 12).Value = 19342.8
$ws.Cells.Item(32, 13).Value = -9038.529
$ws.Cells.Item(32, 14).Value = -19916.8
$ws.Cells.Item(61, 8).Value = 3358.5
$ws.Cells.Item(61, 9).Value = 3561.8
$ws.Cells.Item(61, 10).Value = 2654.7693
$ws.Cells.Item(61, 11).Value = 3561.8
$ws.Cells.Item(61, 12).Value = 2654.7693
$ws.Cells.Item(61, 13).Value = -3349.8
$ws.Cells.Item(61, 14).Value = -3078.7693
$ws.Cells.Item(74, 8).Value = 821
$ws.Cells.Item(74, 9).Value = 670.1
$ws.Cells.Item(74, 10).Value = 1122.8
$ws.Cells.Item(74, 11).Value = 670.1
$ws.Cells.Item(74, 12).Value = 1122.8
$ws.Cells.Item(74, 13).Value = 203.9
$ws.Cells.Item(74, 14).Value = -2870.8
$ws.Cells.Item(77, 8).Value = 821
$ws.Cells.Item(77, 9).Value = 670.1
$ws.Cells.Item(77, 10).Value = 1122.8
$ws.Cells.Item(77, 11).Value = 3350.5
$ws.Cells.Item(77, 12).Value = 5614
$ws.Cells.Item(77, 13).Value = 1017.5
$ws.Cells.Item(77, 14).Value = -14350
$ws.Cells.Item(102, 8).Value = 47620540
$ws.Cells.Item(102, 9).Value = 76924520
$ws.Cells.Item(102, 10).Value = 1565.25
$ws.Cells.Item(102, 11).Value = 76924520
$ws.Cells.Item(102, 12).Value = 1565.25
$ws.Cells.Item(102, 13).Value = -76922898
$ws.Cells.Item(102, 14).Value = -4809.25
$ws.Cells.Item(136, 8).Value = 3358.5
$ws.Cells.Item(136, 9).Value = 3561.8
$ws.Cells.Item(136, 10).Value = 2654.7693
$ws.Cells.Item(136, 11).Value = 10685.4
$ws.Cells.Item(136, 12).Value = 7964.3079
$ws.Cells.Item(136, 13).Value = -8135.400000000001
$ws.Cells.Item(136, 14).Value = -13064.3079

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 705.7692
$ws.Cells.Item(107, 9).Value = 536.2
$ws.Cells.Item(107, 10).Value = 1271
$ws.Cells.Item(107, 11).Value = 536.2
$ws.Cells.Item(107, 12).Value = 1271
$ws.Cells.Item(107, 13).Value = 1383.8
$ws.Cells.Item(107, 14).Value = -5111
$ws.Cells.Item(134, 8).Value = 12365833
$ws.Cells.Item(134, 9).Value = 14515787
$ws.Cells.Item(134, 11).Value = 43547361
$ws.Cells.Item(134, 13).Value = -43544826

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4449.6353
$ws.Cells.Item(31, 9).Value = 955.34546
$ws.Cells.Item(31, 10).Value = 10855.833
$ws.Cells.Item(31, 11).Value = 955.34546
$ws.Cells.Item(31, 12).Value = 10855.833
$ws.Cells.Item(31, 13).Value = -660.34546
$ws.Cells.Item(31, 14).Value = -11445.833
$ws.Cells.Item(34, 8).Value = 4449.6353
$ws.Cells.Item(34, 9).Value = 955.34546
$ws.Cells.Item(34, 10).Value = 10855.833
$ws.Cells.Item(34, 11).Value = 955.34546
$ws.Cells.Item(34, 12).Value = 10855.833
$ws.Cells.Item(34, 13).Value = -753.34546
$ws.Cells.Item(34, 14).Value = -11259.833
$ws.Cells.Item(58, 8).Value = 3348927.5
$ws.Cells.Item(58, 9).Value = 3996985.2
$ws.Cells.Item(58, 10).Value = 16059
$ws.Cells.Item(58, 11).Value = 3996985.2
$ws.Cells.Item(58, 12).Value = 16059
$ws.Cells.Item(58, 13).Value = -3996782.2
$ws.Cells.Item(58, 14).Value = -16465
$ws.Cells.Item(132, 8).Value = 6292438.5
$ws.Cells.Item(132, 9).Value = 7937197
$ws.Cells.Item(132, 11).Value = 23811591
$ws.Cells.Item(132, 13).Value = -23809061
$ws.Cells.Item(134, 8).Value = 7813746
$ws.Cells.Item(134, 9).Value = 7354119
$ws.Cells.Item(134, 11).Value = 22062357
$ws.Cells.Item(134, 13).Value = -22059822
$ws.Cells.Item(136, 8).Value = 3348927.5
$ws.Cells.Item(136, 9).Value = 3996985.2
$ws.Cells.Item(136, 10).Value = 16059
$ws.Cells.Item(136, 11).Value = 11990955.6
$ws.Cells.Item(136, 12).Value = 48177
$ws.Cells.Item(136, 13).Value = -11988405.6
$ws.Cells.Item(136, 14).Value = -53277

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 727
$ws.Cells.Item(117, 9).Value = 400
$ws.Cells.Item(117, 10).Value = 792.4
$ws.Cells.Item(117, 11).Value = 1200
$ws.Cells.Item(117, 12).Value = 2377.2
$ws.Cells.Item(117, 13).Value = 2242
$ws.Cells.Item(117, 14).Value = -9261.200000000001
$ws.Cells.Item(129, 8).Value = 1308
$ws.Cells.Item(129, 9).Value = 470
$ws.Cells.Item(129, 10).Value = 1601.3
$ws.Cells.Item(129, 11).Value = 1410
$ws.Cells.Item(129, 12).Value = 4803.9
$ws.Cells.Item(129, 13).Value = 3590
$ws.Cells.Item(129, 14).Value = -14803.9
$ws.Cells.Item(132, 8).Value = 1754.25
$ws.Cells.Item(132, 9).Value = 639.2632
$ws.Cells.Item(132, 11).Value = 5753.3688
$ws.Cells.Item(132, 13).Value = -3223.3688

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 55556496
$ws.Cells.Item(97, 9).Value = 100000850
$ws.Cells.Item(97, 10).Value = 1051.375
$ws.Cells.Item(97, 11).Value = 100000850
$ws.Cells.Item(97, 12).Value = 1051.375
$ws.Cells.Item(97, 13).Value = -100000354
$ws.Cells.Item(97, 14).Value = -2043.375
$ws.Cells.Item(132, 8).Value = 15643528
$ws.Cells.Item(132, 9).Value = 19629002
$ws.Cells.Item(132, 11).Value = 58887006
$ws.Cells.Item(132, 13).Value = -58884476

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1242.5
$ws.Cells.Item(100, 9).Value = 1155.5555
$ws.Cells.Item(100, 10).Value = 1354.2858
$ws.Cells.Item(100, 11).Value = 1155.5555
$ws.Cells.Item(100, 12).Value = 1354.2858
$ws.Cells.Item(100, 13).Value = -614.5554999999999
$ws.Cells.Item(100, 14).Value = -2436.2858

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 10003070
$ws.Cells.Item(81, 9).Value = 33333900
$ws.Cells.Item(81, 10).Value = 4142.857
$ws.Cells.Item(81, 11).Value = 66667800
$ws.Cells.Item(81, 12).Value = 8285.714
$ws.Cells.Item(81, 13).Value = -66666739
$ws.Cells.Item(81, 14).Value = -10407.714
$ws.Cells.Item(84, 8).Value = 10003070
$ws.Cells.Item(84, 9).Value = 33333900
$ws.Cells.Item(84, 10).Value = 4142.857
$ws.Cells.Item(84, 11).Value = 333339000
$ws.Cells.Item(84, 12).Value = 41428.57
$ws.Cells.Item(84, 13).Value = -333333696
$ws.Cells.Item(84, 14).Value = -52036.57
$ws.Cells.Item(132, 8).Value = 246814.86
$ws.Cells.Item(132, 9).Value = 21640.143
$ws.Cells.Item(132, 10).Value = 1166278.2
$ws.Cells.Item(132, 11).Value = 64920.429
$ws.Cells.Item(132, 12).Value = 3498834.6
$ws.Cells.Item(132, 13).Value = -62390.429
$ws.Cells.Item(132, 14).Value = -3503894.6
$ws.Cells.Item(136, 8).Value = 699683.9
$ws.Cells.Item(136, 9).Value = 1033667.94
$ws.Cells.Item(136, 10).Value = 1353.5454
$ws.Cells.Item(136, 11).Value = 3101003.82
$ws.Cells.Item(136, 12).Value = 4060.6362
$ws.Cells.Item(136, 13).Value = -3098453.82
$ws.Cells.Item(136, 14).Value = -9160.636200000001

Write-Host "Applied all changes"